# Natmi following Dr Hou advice
# Update the Ptn-Ptprs LR-pair sheet with recomputed NATMI statistics
# (ligand/receptor-expressing cell counts changed from 1 to 3, which in
# turn changes the derived total-expression and specificity columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2 = @{
        E = 3;  G = 3.270036666666666;  H = 9.81011;
        I = 0.359406393324744;          J = 0.3594063933247441;
        K = 3;  M = 2.618716333333334;  N = 7.856149000000001;
        O = 0.07115908183301342;        P = 0.07115908183301341;
        Q = 8.563298429598889;          R = 77.06968586639002;
        S = 0.02557502895390367;        T = 0.02557502895390367;
    }
    3 = @{
        E = 3;  G = 3.270036666666666;  H = 9.81011;
        I = 0.359406393324744;          J = 0.3594063933247441;
        K = 3;  M = 15.503283;          N = 46.509849;
        O = 0.4212748702999519;         P = 0.4212748702999519;
        Q = 50.69630386371;             R = 456.26673477339;
        S = 0.151408881732855;          T = 0.151408881732855;
    }
    4 = @{
        E = 3;  G = 3.270036666666666;  H = 9.81011;
        I = 0.359406393324744;          J = 0.3594063933247441;
        K = 3;  M = 18.67887366666666;  N = 56.036621;
        O = 0.5075660478670347;         P = 0.5075660478670347;
        Q = 61.08060178203444;          R = 549.72541603831;
        S = 0.1824224826379853;         T = 0.1824224826379854;
    }
    5 = @{
        E = 3;  G = 5.828401;           H = 17.485203;
        I = 0.6405936066752559;         J = 0.640593606675256;
        K = 3;  M = 2.618716333333334;  N = 7.856149000000001;
        O = 0.07115908183301342;        P = 0.07115908183301341;
        Q = 15.26292889591634;          R = 137.366360063247;
        S = 0.04558405287910974;        T = 0.04558405287910974;
    }
    6 = @{
        E = 3;  G = 5.828401;           H = 17.485203;
        I = 0.6405936066752559;         J = 0.640593606675256;
        K = 3;  M = 15.503283;          N = 46.509849;
        O = 0.4212748702999519;         P = 0.4212748702999519;
        Q = 90.35935014048302;          R = 813.2341512643471;
        S = 0.2698659885670968;         T = 0.2698659885670969;
    }
    7 = @{
        E = 3;  G = 5.828401;           H = 17.485203;
        I = 0.6405936066752559;         J = 0.640593606675256;
        K = 3;  M = 18.67887366666666;  N = 56.036621;
        O = 0.5075660478670347;         P = 0.5075660478670347;
        Q = 108.8679659576737;          R = 979.811693619063;
        S = 0.3251435652290494;         T = 0.3251435652290494;
    }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
